# This script reproduces the authoritative edit to the 'Avverkningsanmälningar' sheet:
#  1. Every data row's 'Förändrad' date (column C) advances from 46070 to 46072.
#  2. Within rows 12-16, the last entry (A 1793-2025) moves up to become the first
#     (row 12), shifting the other four rows down by one.
#  3. Within rows 25-31, the first entry (A 63378-2025) moves down to become the
#     last (row 31), shifting the other six rows up by one.
# Row contents are written explicitly (rather than relying on Cut/Insert semantics)
# so the result is unambiguous; cells whose value would be unchanged are skipped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update 'Förändrad' (column C) for every data row 2-39 to the new date serial 46072
for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 3).Value = 46072
}

# 2) Re-home the rows whose content shifted position (blocks 12-16 and 25-31).
#    (Column C is intentionally skipped here -- already set to 46072 above.)

# row 12 <- original row 16 (A='A 1793-2025')
$ws.Range('A12').Value = 'A 1793-2025'
$ws.Range('B12').Value = 45671
$ws.Range('G12').Value = 1.7
$ws.Range('H12').Value = 1
$ws.Range('J12').Value = 0
$ws.Range('O12').Value = 0
$ws.Range('R12').Value = 'Mistel'
$ws.Range('S12').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/artfynd/A 1793-2025 artfynd.xlsx", "A 1793-2025")'
$ws.Range('T12').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/kartor/A 1793-2025 karta.png", "A 1793-2025")'
$ws.Range('V12').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/klagomål/A 1793-2025 FSC-klagomål.docx", "A 1793-2025")'
$ws.Range('W12').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/klagomålsmail/A 1793-2025 FSC-klagomål mail.docx", "A 1793-2025")'
$ws.Range('X12').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/tillsyn/A 1793-2025 tillsynsbegäran.docx", "A 1793-2025")'
$ws.Range('Y12').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/tillsynsmail/A 1793-2025 tillsynsbegäran mail.docx", "A 1793-2025")'

# row 13 <- original row 12 (A='A 13794-2023')
$ws.Range('A13').Value = 'A 13794-2023'
$ws.Range('B13').Value = 45007
$ws.Range('G13').Value = 4.8
$ws.Range('I13').Value = 0
$ws.Range('J13').Value = 1
$ws.Range('O13').Value = 1
$ws.Range('R13').Value = 'Tallticka'
$ws.Range('S13').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/artfynd/A 13794-2023 artfynd.xlsx", "A 13794-2023")'
$ws.Range('T13').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/kartor/A 13794-2023 karta.png", "A 13794-2023")'
$ws.Range('V13').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/klagomål/A 13794-2023 FSC-klagomål.docx", "A 13794-2023")'
$ws.Range('W13').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/klagomålsmail/A 13794-2023 FSC-klagomål mail.docx", "A 13794-2023")'
$ws.Range('X13').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/tillsyn/A 13794-2023 tillsynsbegäran.docx", "A 13794-2023")'
$ws.Range('Y13').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/tillsynsmail/A 13794-2023 tillsynsbegäran mail.docx", "A 13794-2023")'

# row 14 <- original row 13 (A='A 49128-2025')
$ws.Range('A14').Value = 'A 49128-2025'
$ws.Range('B14').Value = 45937.8944675926
$ws.Range('G14').Value = 3.6
$ws.Range('H14').Value = 0
$ws.Range('I14').Value = 1
$ws.Range('R14').Value = 'Rödgul trumpetsvamp'
$ws.Range('S14').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/artfynd/A 49128-2025 artfynd.xlsx", "A 49128-2025")'
$ws.Range('T14').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/kartor/A 49128-2025 karta.png", "A 49128-2025")'
$ws.Range('V14').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/klagomål/A 49128-2025 FSC-klagomål.docx", "A 49128-2025")'
$ws.Range('W14').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/klagomålsmail/A 49128-2025 FSC-klagomål mail.docx", "A 49128-2025")'
$ws.Range('X14').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/tillsyn/A 49128-2025 tillsynsbegäran.docx", "A 49128-2025")'
$ws.Range('Y14').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/tillsynsmail/A 49128-2025 tillsynsbegäran mail.docx", "A 49128-2025")'

# row 15 <- original row 14 (A='A 35953-2023')
$ws.Range('A15').Value = 'A 35953-2023'
$ws.Range('B15').Value = 45148
$ws.Range('G15').Value = 2.1
$ws.Range('I15').Value = 0
$ws.Range('R15').Value = 'Vanlig groda'
$ws.Range('S15').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/artfynd/A 35953-2023 artfynd.xlsx", "A 35953-2023")'
$ws.Range('T15').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/kartor/A 35953-2023 karta.png", "A 35953-2023")'
$ws.Range('V15').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/klagomål/A 35953-2023 FSC-klagomål.docx", "A 35953-2023")'
$ws.Range('W15').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/klagomålsmail/A 35953-2023 FSC-klagomål mail.docx", "A 35953-2023")'
$ws.Range('X15').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/tillsyn/A 35953-2023 tillsynsbegäran.docx", "A 35953-2023")'
$ws.Range('Y15').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/tillsynsmail/A 35953-2023 tillsynsbegäran mail.docx", "A 35953-2023")'

# row 16 <- original row 15 (A='A 52545-2024')
$ws.Range('A16').Value = 'A 52545-2024'
$ws.Range('B16').Value = 45609.63619212963
$ws.Range('G16').Value = 2.8
$ws.Range('I16').Value = 1
$ws.Range('R16').Value = 'Grön sköldmossa'
$ws.Range('S16').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/artfynd/A 52545-2024 artfynd.xlsx", "A 52545-2024")'
$ws.Range('T16').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/kartor/A 52545-2024 karta.png", "A 52545-2024")'
$ws.Range('V16').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/klagomål/A 52545-2024 FSC-klagomål.docx", "A 52545-2024")'
$ws.Range('W16').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/klagomålsmail/A 52545-2024 FSC-klagomål mail.docx", "A 52545-2024")'
$ws.Range('X16').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/tillsyn/A 52545-2024 tillsynsbegäran.docx", "A 52545-2024")'
$ws.Range('Y16').Formula = '=HYPERLINK("https://klasma.github.io/Logging_0125/tillsynsmail/A 52545-2024 tillsynsbegäran mail.docx", "A 52545-2024")'

# row 25 <- original row 26 (A='A 34073-2025')
$ws.Range('A25').Value = 'A 34073-2025'
$ws.Range('B25').Value = 45845
$ws.Range('G25').Value = 2.9

# row 26 <- original row 27 (A='A 35047-2025')
$ws.Range('A26').Value = 'A 35047-2025'
$ws.Range('B26').Value = 45852.35094907408
$ws.Range('G26').Value = 4.2

# row 27 <- original row 28 (A='A 37245-2025')
$ws.Range('A27').Value = 'A 37245-2025'
$ws.Range('B27').Value = 45876.39396990741
$ws.Range('G27').Value = 1.1

# row 28 <- original row 29 (A='A 2533-2026')
$ws.Range('A28').Value = 'A 2533-2026'
$ws.Range('B28').Value = 46037.44622685185
$ws.Range('G28').Value = 1.4

# row 29 <- original row 30 (A='A 2535-2026')
$ws.Range('A29').Value = 'A 2535-2026'
$ws.Range('B29').Value = 46037.4490625
$ws.Range('G29').Value = 1.6

# row 30 <- original row 31 (A='A 55724-2024')
$ws.Range('A30').Value = 'A 55724-2024'
$ws.Range('B30').Value = 45622.81020833334
$ws.Range('G30').Value = 0.5

# row 31 <- original row 25 (A='A 63378-2025')
$ws.Range('A31').Value = 'A 63378-2025'
$ws.Range('B31').Value = 46010.72528935185
$ws.Range('G31').Value = 4.8
